# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Changes applied:
#  1. "Periodo Mora" updated from 2507 to 2508 for every worker row
#     (E16, E17, E18 all shared the same value).
#  2. "Valor Mora" for the first worker (row 16, column G) updated
#     from 4348000 to 1870000.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update Periodo Mora (column E) from 2507 -> 2508 for all data rows
$ws.Range("E16").Value = "2508"
$ws.Range("E17").Value = "2508"
$ws.Range("E18").Value = "2508"

# 2. Update Valor Mora for first worker (G16) from 4348000 -> 1870000
$ws.Range("G16").Value = 1870000
